$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

# --- F2: updated SQL query text (column label + two formula tweaks) ---
$query = @'
SELECT M.[AgentID] as 'Agent ID', A.[AgentName] as 'Agent Name' ,A.TeamName as 'Team Name',A.SupervisorName 'Supervisor Name' ,SUM([TotalInteraction]) [Total Interaction] ,SUM([TotalVoice]) [Total Voice],SUM(TotalChat) [Total Chat],
SUM([TotalAudioIP]) [Total Audio IP],SUM([TotalVideoIP]) [Total Video IP],SUM(TotalSM) [Total SM],SUM(TotalSMS) [Total SMS],SUM(isnull(TotalFax,0)) [Total Fax],SUM(TotalEmail) [Total Email],  
[dbo].[SECONDSTOhhmmss](SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms,M.TotalStaffedTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalStaffedTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar,DATEADD(ms,M.TotalStaffedTime* 1000, 0), 108)))) as [Total Staffed Time] ,  	  
[dbo].[SECONDSTOhhmmss](SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalInteractionTime* 1000, 0), 108)) * 3600 +DATEPART(mi,CONVERT(varchar, DATEADD(ms, M.TotalInteractionTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar,DATEADD(ms,M.TotalInteractionTime* 1000, 0), 108))))as [Total Interaction Time] ,   
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalInteractionTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalInteractionTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalInteractionTime* 1000, 0), 108)))/nullif(SUM([TotalInteraction]),0),0)) as [Avg Interaction Time],
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalACWTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalACWTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalACWTime* 1000, 0), 108)))) as [Total ACW Time], 
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalAuxTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalAuxTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalAuxTime* 1000, 0), 108)))) as [Total Aux Time],
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalVoiceTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalVoiceTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalVoiceTime* 1000, 0), 108)))) as [Total Voice Time],	
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,[TotalVoiceTime]) * 3600 +DATEPART(mi, [TotalVoiceTime]) * 60 + DATEPART(ss,[TotalVoiceTime]))/nullif(SUM(TotalVoice),0),0))as [Avg Voice Time] ,  
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)))) as [Total Chat Time],  
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalChatTime* 1000, 0), 108)))/nullif(SUM([TotalChat]),0),0)) as [Avg Chat Time],
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalSMTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalSMTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalSMTime* 1000, 0), 108)))) as [Total SM Time],
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,[TotalSMTime]) * 3600 +DATEPART(mi, [TotalSMTime]) * 60 + DATEPART(ss,[TotalSMTime]))/nullif(SUM(TotalSM),0),0))as [Avg SM Time] ,  
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)))) as [Total SMS Time],
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalSMSTime* 1000, 0), 108)))/nullif(SUM([TotalSMS]),0),0)) as [Avg SMS Time],
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)))) as [Total Fax Time],	
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalFaxTime* 1000, 0), 108)))/nullif(SUM([TotalFax]),0),0)) as [Avg Fax Time],
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalEmailTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalEmailTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalEmailTime* 1000, 0), 108)))) as [Total Email Time],
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,[TotalEmailTime]) * 3600 +DATEPART(mi, [TotalEmailTime]) * 60 + DATEPART(ss,[TotalEmailTime]))/nullif(SUM(TotalEmail),0),0))as [Avg Email Time] ,  
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)))) as [Total Audio IP Time],
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalAudioIPTime* 1000, 0), 108)))/nullif(SUM([TotalAudioIP]),0),0)) as [Avg Audio IP Time],  
[dbo].[SECONDSTOhhmmss](SUM(DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)))) as [Total Video IP Time],
[dbo].[SECONDSTOhhmmss](ISNULL(SUM ( DATEPART(hh,CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)) * 3600 +DATEPART(mi, CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)) * 60 + DATEPART(ss,CONVERT(varchar, DATEADD(ms, M.TotalVideoIPTime* 1000, 0), 108)))/nullif(SUM([TotalVideoIP]),0),0)) as [Avg Video IP Time],  
SUM(TotalExtIn)as [Total Ext In] ,SUM(TotalExtOut) as [Total Ext Out] ,SUM(TotalTransferIn) as [Total Transfer In] ,SUM(TotalTransferOut) as [Total Transfer Out] ,SUM(TotalConferenceIn) as [Total Conference In] ,  SUM(TotalConferenceOut) as [Total Conference Out] 
From (Select TAB2.AgentID as AgentID, Isnull(FirstName,'')+' '+isnull(LastName,'') AgentName,TAB2.StationID,TAB2.LoginDateTime as LoginDateTime,TAB2.LogoutDateTime as LogoutDateTime,
		TAB2.SkillList, Isnull(TAB2.TotalStaffedTime,0) as TotalStaffedTime,Skills.SkillLists AS SkillNameList, 
		Isnull(SUM(TotalACWTime),0) TotalACWTime, Isnull(isnull(TAB2.TotalStaffedTime,0) - (isnull(SUM(TotalInteractionTime),0) + isnull(SUM(TotalACWTime),0)),0) TotalAuxTime,
		Isnull(SUM(TotalInteraction),0) TotalInteraction,Isnull(SUM(TotalInteractionTime)+SUM(TotalInteractionHoldTime),0) TotalInteractionTime,
		Isnull(SUM(TotalVoice),0) TotalVoice,Isnull(SUM(TotalVoiceTime)+SUM(TotalVoiceHoldTime),0) TotalVoiceTime,
		Isnull(SUM(TotalEmail),0) TotalEmail,Isnull(SUM(TotalEmailTime)+SUM(TotalEmailHoldTime),0) TotalEmailTime,
		Isnull(SUM(TotalChat),0)TotalChat,Isnull(SUM(TotalChatTime)+SUM(TotalChatHoldTime),0) TotalChatTime,
		Isnull(SUM(TotalSMS),0)TotalSMS,Isnull(SUM(TotalSMSTime)+SUM(TotalSMSHoldTime),0) TotalSMSTime,
		Isnull(SUM(TotalFax),0)TotalFax,Isnull(SUM(TotalFaxTime)+SUM(TotalFaxHoldTime),0) TotalFaxTime,
		Isnull(SUM(TotalSM),0)TotalSM,Isnull(SUM(TotalSMTime)+SUM(TotalSMHoldTime),0) TotalSMTime,
		Isnull(SUM(TotalAudioIP),0)TotalAudioIP,Isnull(SUM(TotalAudioIPTime)+SUM(TotalAudioIPHoldTime),0) TotalAudioIPTime,
		Isnull(SUM(TotalVideoIP),0)TotalVideoIP,Isnull(SUM(TotalVideoIPTime)+SUM(TotalVideoIPHoldTime),0) TotalVideoIPTime,
		Isnull(SUM(TotalExtIn),0)TotalExtIn,
		Isnull(SUM(TotalExtOut),0)TotalExtOut,Isnull(SUM(TotalTranIn),0)TotalTransferIn,Isnull(SUM(TotalTranOut),0) TotalTransferOut,
		Isnull(SUM(TotalConfIn),0)TotalConferenceIn,Isnull(SUM(TotalConfOut),0)TotalConferenceOut FROM (
		Select TMAC_Interactions.AgentID, ClosedDateTime,
		TotalInteraction = CASE WHEN Channel IN ('Voice','Chat','TextChat','FBPost','FBPrivate','Email','AudioChat','VideoChat','SMS','FAX') THEN  1 ELSE 0 END,
		TotalInteractionTime = CASE WHEN Channel IN ('Voice','Chat','TextChat','FBPost','FBPrivate','Email','AudioChat','VideoChat','SMS','FAX') THEN  ActiveTime ELSE 0 END,
		TotalInteractionHoldTime = CASE WHEN Channel IN ('Voice','Chat','TextChat','FBPost','FBPrivate','Email','AudioChat','VideoChat','SMS','FAX') THEN  HoldTime ELSE 0 END,

		TotalACWTime  = CASE WHEN Channel IN ('Voice','Chat','TextChat','FBPost','FBPrivate','Email','AudioChat','VideoChat','SMS','FAX')  THEN ACWTime ELSE 0 END,

		TotalVoice  = CASE WHEN (Channel='Voice' AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalEmail  = CASE WHEN (Channel='Email' AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalChat  = CASE WHEN (Channel IN ('Chat','TextChat') AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalSMS  = CASE WHEN (Channel='SMS' AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalFax  = CASE WHEN (Channel='FAX' AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalSM  = CASE WHEN (Channel IN ('FBPost','FBPrivate') AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalAudioIP  = CASE WHEN ((SubChannel='Audio' OR Channel='AudioChat') AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,
		TotalVideoIP  = CASE WHEN ((SubChannel='Video' OR Channel='VideoChat') AND IsTransfered<>1 AND IsConferenced<>1) THEN 1 ELSE 0 END,

		TotalVoiceTime  = CASE WHEN Channel='Voice' THEN ActiveTime ELSE 0 END,
		TotalEmailTime  = CASE WHEN Channel='Email' THEN ActiveTime ELSE 0 END,
		TotalChatTime  = CASE WHEN Channel IN ('Chat','TextChat') THEN ActiveTime ELSE 0 END,
		TotalSMSTime  = CASE WHEN Channel='SMS' THEN ActiveTime ELSE 0 END,
		TotalFaxTime  = CASE WHEN Channel='FAX' THEN ActiveTime ELSE 0 END,
		TotalSMTime  = CASE WHEN Channel IN ('FBPost','FBPrivate') THEN ActiveTime ELSE 0 END,
		TotalAudioIPTime  = CASE WHEN (SubChannel='Audio' OR Channel='AudioChat') THEN ActiveTime ELSE 0 END,
		TotalVideoIPTime  = CASE WHEN (SubChannel='Video' OR Channel='VideoChat') THEN ActiveTime ELSE 0 END,

		TotalVoiceHoldTime  = CASE WHEN Channel='Voice' THEN HoldTime ELSE 0 END,
		TotalEmailHoldTime  = CASE WHEN Channel='Email' THEN HoldTime ELSE 0 END,
		TotalChatHoldTime  = CASE WHEN Channel IN ('Chat','TextChat') THEN HoldTime ELSE 0 END,
		TotalSMSHoldTime  = CASE WHEN Channel='SMS' THEN HoldTime ELSE 0 END,
		TotalFaxHoldTime  = CASE WHEN Channel='FAX' THEN HoldTime ELSE 0 END,
		TotalSMHoldTime  = CASE WHEN Channel IN ('FBPost','FBPrivate') THEN HoldTime ELSE 0 END,
		TotalAudioIPHoldTime  = CASE WHEN (SubChannel='Audio' OR Channel='AudioChat') THEN HoldTime ELSE 0 END,
		TotalVideoIPHoldTime  = CASE WHEN (SubChannel='Video' OR Channel='VideoChat') THEN HoldTime ELSE 0 END,
		TotalExtIn  = CASE WHEN (CallType='2' AND Direction='In' AND Channel='Voice') THEN 1 ELSE 0 END,
		TotalExtOut  = CASE WHEN (CallType='2' AND Direction='Out' AND Channel='Voice') THEN 1 ELSE 0 END,
		TotalTranIn  = CASE WHEN IsTransfered='1' THEN 1 ELSE 0 END,
		TotalTranOut  = CASE WHEN IsTranferedTo='1' THEN 1 ELSE 0 END,
		TotalConfIn  = CASE WHEN IsConferenced='1' THEN 1 ELSE 0 END,
		TotalConfOut = CASE WHEN IsConferencedTo='1' THEN 1 ELSE 0 END,LoginInstanceID

		From TMAC_Interactions with(nolock) 
		) AS TAB1 
		RIGHT JOIN
		(SELECT  AgentID,LoginDateTime,LogoutDateTime,A.LoginInstanceID,
		SUM(DATEDIFF(SECOND,CONVERT(datetime, STUFF(STUFF(STUFF(LoginDateTime, 9, 0, ' '), 12, 0, ':'), 15, 0, ':')),
		CONVERT(datetime,STUFF(STUFF(STUFF(LogoutDateTime, 9, 0, ' '), 12, 0, ':'), 15, 0, ':')))) AS TotalStaffedTime,StationID,SkillList
		FROM (SELECT DISTINCT AgentID,LoginDate+LoginTime AS LoginDateTime ,LogoutDate+LogoutTime AS LogoutDateTime,LoginInstanceID,StationID,SkillList
		FROM AGT_Agent_TimeTrack A with(nolock) WHERE LogoutDate+LogoutTime>='ReportBeforeDate' AND LogoutDate+LogoutTime<='ReportAfterDate' 
		)A Group by AgentID,StationID,SkillList,LoginDateTime,LogoutDateTime,LoginInstanceID) AS TAB2 
		ON TAB2.AgentID=TAB1.AgentId AND TAB2.LoginInstanceID = TAB1.LoginInstanceID AND TAB1.ClosedDateTime between LoginDateTime AND LogoutDateTime 
		INNER JOIN AGT_Agent A with(nolock) ON A.AvayaLoginID = TAB2.AgentID
		CROSS APPLY 
		(SELECT STUFF((SELECT ','+ SkillName from (
		SELECT SkillName  FROM TMAC_Skills
		WHERE SkillExtension in (Select Id from BreakStringIntoRows (TAB2.SkillList)) 
		union
		SELECT ID as SkillName FROM BreakStringIntoRows (TAB2.SkillList) WHERE ID NOT IN (SELECT SkillExtension from TMAC_Skills)) Skill order by SkillName desc
		FOR XML PATH('')), 1, 1, '' ) SkillLists
		) AS Skills
		Group by TAB2.AgentID,FirstName,LastName,TAB2.StationID,TAB2.LoginDateTime,TAB2.LogoutDateTime,TAB2.SkillList,
		TAB2.TotalStaffedTime,Skills.SkillLists 
		--ORDER BY LoginDateTime ASC
		) M  INNER JOIN fn_AgentHierarchy('na',1,1) A on A.AgentId=M.AgentID
		GROUP BY M.[AgentID],A.[AgentName],A.TeamName,A.SupervisorName
'@
$ws.Range("F2").Value = $query

# --- Date cells D2/E2: typed with a leading apostrophe (text/quote-prefix entry),
# matching the existing quotePrefix cell style already used on this sheet. ---
$ws.Range("D2").Formula = "'17-04-2020 00:00:00"
$ws.Range("E2").Formula = "'20-04-2020 00:00:00"

# F2 alignment: drop the explicit "top" vertical alignment, keep wrap text
$ws.Range("F2").VerticalAlignment = -4107

# Column widths: D/E share a width, F gets a bit wider
$ws.Columns.Item(4).ColumnWidth = 17.0
$ws.Columns.Item(5).ColumnWidth = 17.0
$ws.Columns.Item(6).ColumnWidth = 19.67

# Selection moves from F2 to E2
[void]$ws.Range("E2").Select()
